$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" '30.436.99'
$ws.Range("E2").Value = '  +0.42%  '

Set-TextValue "D3" '1.856.04'
$ws.Range("E3").Value = '  -0.24%  '

Set-TextValue "D4" '0.9996'
$ws.Range("E4").Value = '  +0.01%  '

Set-TextValue "D5" '233.69'
$ws.Range("E5").Value = '  +0.13%  '

Set-TextValue "D6" '0.9999'
$ws.Range("E6").Value = '  +0.03%  '

Set-TextValue "D7" '0.4682'
$ws.Range("E7").Value = '  -1.57%  '

Set-TextValue "D8" '0.2750'
$ws.Range("E8").Value = '  -0.38%  '

Set-TextValue "D9" '0.06320'
$ws.Range("E9").Value = '  -2.00%  '

Set-TextValue "D10" '1.864.42'
$ws.Range("E10").Value = '  +0.32%  '

Set-TextValue "D11" '17.14'
$ws.Range("E11").Value = '  +6.10%  '

Set-TextValue "D12" '0.07460'
$ws.Range("E12").Value = '  +0.61%  '

Set-TextValue "D13" '4.949'
$ws.Range("E13").Value = '  -1.19%  '

Set-TextValue "D14" '84.20'
$ws.Range("E14").Value = '  -1.90%  '

Set-TextValue "D15" '0.6272'
$ws.Range("E15").Value = '  -1.19%  '

Set-TextValue "D16" '30.378.40'
$ws.Range("E16").Value = '  +0.33%  '

$ws.Range("E17").Value = '  +0.11%  '

Set-TextValue "D18" '229.53'
$ws.Range("E18").Value = '  +0.37%  '

Set-TextValue "D19" '12.53'
$ws.Range("E19").Value = '  -2.38%  '

Set-TextValue "D20" '0.000007326'
$ws.Range("E20").Value = '  -0.91%  '

$ws.Range("E21").Value = '  -0.06%  '

Set-TextValue "D22" '4.937'
$ws.Range("E22").Value = '  -3.56%  '

Set-TextValue "D23" '5.914'
$ws.Range("E23").Value = '  -2.06%  '

Set-TextValue "D24" '166.84'

Set-TextValue "D25" '9.207'
$ws.Range("E25").Value = '  -0.91%  '

Set-TextValue "D26" '17.91'
$ws.Range("E26").Value = '  -0.02%  '

Set-TextValue "D27" '1.883'
$ws.Range("E27").Value = '  +1.01%  '

$ws.Range("E28").Value = '  -0.97%  '

$ws.Range("E29").Value = '  -0.64%  '

Set-TextValue "D30" '4.098'
$ws.Range("E30").Value = '  -3.50%  '

Set-TextValue "D31" '3.831'
$ws.Range("E31").Value = '  -2.24%  '

$ws.Range("E32").Value = '  -0.26%  '

Set-TextValue "D33" '1.143'
$ws.Range("E33").Value = '  -0.71%  '

Set-TextValue "D34" '0.7057'
$ws.Range("E34").Value = '  -3.18%  '

Set-TextValue "D35" '2.708'
$ws.Range("E35").Value = '  +0.72%  '

Set-TextValue "D36" '0.01915'
$ws.Range("E36").Value = '  -2.40%  '

Set-TextValue "D37" '2.672'
$ws.Range("E37").Value = '  +1.60%  '

Set-TextValue "D38" '0.8689'
$ws.Range("E38").Value = '  -4.29%  '

Set-TextValue "D39" '1.953'
$ws.Range("E39").Value = '  -1.86%  '

Set-TextValue "D40" '105.84'
$ws.Range("E40").Value = '  -0.14%  '

Set-TextValue "D41" '0.9999'
$ws.Range("E41").Value = '  +0.05%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue "D42" '0.4060'
$ws.Range("E42").Value = '  -1.61%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D43" '5.515'
$ws.Range("E43").Value = '  -1.02%  '

Set-TextValue "D44" '7.138'
$ws.Range("E44").Value = '  +0.81%  '

Set-TextValue "D45" '61.56'
$ws.Range("E45").Value = '  +0.38%  '

Set-TextValue "D46" '0.1218'
$ws.Range("E46").Value = '  +0.64%  '

Set-TextValue "D47" '33.47'
$ws.Range("E47").Value = '  +1.37%  '

Set-TextValue "D48" '8.603'
$ws.Range("E48").Value = '  -2.68%  '

Set-TextValue "D49" '0.05556'
$ws.Range("E49").Value = '  -0.97%  '

Set-TextValue "D50" '1.367'
$ws.Range("E50").Value = '  -2.59%  '

Set-TextValue "D51" '0.3682'
$ws.Range("E51").Value = '  -1.38%  '

